$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update Latest Handback DateTime for the ca05378c row (row 2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-11-01 14:37:11"

# "de-de" sheet: update Latest Handback DateTime for the ca05378c row (row 2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-11-01 14:37:28"
